$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number must be forced to Text format
# first, otherwise Excel auto-converts the typed string into a numeric value
# (losing the original text formatting, e.g. "1.00" -> 1).
$ws.Range("D2").Value = "51.053.16"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "2.958.88"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.41"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.48"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("E7").Value = "  +1.98%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.44"
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0851"
$ws.Range("E12").Value = "  +2.31%  "
$ws.Range("D13").Value = "3.427.89"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.43"
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("B15").Value = "Uniswap"
$ws.Range("C15").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "12.38"
$ws.Range("E15").Value = "  +74.18%  "
$ws.Range("E16").Value = "  +5.95%  "
$ws.Range("D17").Value = "2.953.72"
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("E18").Value = "  +4.40%  "
$ws.Range("D19").Value = "51.113.32"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.40"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "0.0₃0956"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.35"
$ws.Range("E23").Value = "  +17.71%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.00"
$ws.Range("E24").Value = "  +3.33%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.66"
$ws.Range("E25").Value = "  +2.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.01"
$ws.Range("E26").Value = "  -2.17%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.166"
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.83"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.03"
$ws.Range("E30").Value = "  -8.24%  "
$ws.Range("E31").Value = "  -3.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.52"
$ws.Range("E32").Value = "  +7.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.36"
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "50.91"
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.06"
$ws.Range("E35").Value = "  +2.32%  "
$ws.Range("E36").Value = "  -4.76%  "
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  +8.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.77"
$ws.Range("E39").Value = "  +3.12%  "
$ws.Range("E40").Value = "  +2.21%  "
$ws.Range("E41").Value = "  +2.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.50"
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.31"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("E44").Value = "  +12.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.67"
$ws.Range("E45").Value = "  +3.13%  "
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").Value = "2.046.18"
$ws.Range("E47").Value = "  +2.24%  "
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.258"
$ws.Range("E49").Value = "  -4.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0323"
$ws.Range("E50").Value = "  -4.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.34"
$ws.Range("E51").Value = "  +7.07%  "
